$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the player's name
$ws.Name = "Marcus Stoinis"

# New header row: a "matchNo" column is inserted before the old column A,
# shifting every other header one column to the right (A:L -> B:M).
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Full set of innings rows (the sheet grows from a single row of data to seven).
$data = @(
    @("22nd","Delhi Capitals","Marcus Stoinis","c †de Villiers b Patel","22","17","3","0","129.41","Royal Challengers Bangalore","Ahmedabad","April 27","RCB won by 1 run"),
    @("Qualifier","Delhi Capitals","Marcus Stoinis","b Shivam Mavi","18","23","1","0","78.26","Kolkata Knight Riders","Sharjah","October 13","KKR won by 3 wickets (with 1 ball remaining)"),
    @("2nd","Delhi Capitals","Marcus Stoinis","c Curran b Thakur","14","9","3","0","155.55","Chennai Super Kings","Wankhede","April 10","Capitals won by 7 wickets (with 8 balls remaining)"),
    @("11th","Delhi Capitals","Marcus Stoinis","","27","13","3","1","207.69","Punjab Kings","Wankhede","April 18","Capitals won by 6 wickets (with 10 balls remaining)"),
    @("25th","Delhi Capitals","Marcus Stoinis","","6","3","1","0","200.00","Kolkata Knight Riders","Ahmedabad","April 29","Capitals won by 7 wickets (with 21 balls remaining)"),
    @("20th","Delhi Capitals","Marcus Stoinis","","2","2","0","0","100.00","Sunrisers Hyderabad","Chennai","April 25","Match tied (Capitals won the one-over eliminator)"),
    @("7th","Delhi Capitals","Marcus Stoinis","c Buttler b Mustafizur Rahman","0","5","0","0","0.00","Rajasthan Royals","Wankhede","April 15","Royals won by 3 wickets (with 2 balls remaining)")
)

# Numeric-looking columns (runs, balls, fours, sixes, sr = columns E..I) are stored
# as text in the source data (note stats like "0.00"/"129.41" need the literal
# string, not a rounded/trimmed number) -- force text formatting before writing
# so Excel doesn't re-interpret them as numbers.
$ws.Range("E2:I8").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($value -ne "") {
            $ws.Cells.Item($r + 2, $c + 1).Value = $value
        }
    }
}

